$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E16").Value = 80
$ws.Range("H16").Value = 40
$ws.Range("K16").Value = 20
$ws.Range("N16").Value = 60
$ws.Range("Q16").Value = 0
$ws.Range("T16").Value = 100
